function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1/2: "Xbits" => "X-bit" text fixes inside the nested "Processing Unit" box ---
$grp75 = Get-ShapeById $s.Shapes 76
$sh82 = Get-ShapeById $grp75.GroupItems 82
$tr82 = $sh82.TextFrame.TextRange
$tr82.Paragraphs(4).Runs(1).Text = "• FPU 32-bit & 64-bit"
$tr82.Paragraphs(5).Runs(1).Text = "• 16-bit Pointers(*)"

# --- 3: Resize/move the "Resources" outer container (id 91) ---
$sh91 = Get-ShapeById $s.Shapes 91
$sh91.Top = 130.0138702392578
$sh91.Height = 248.5040283203125

# --- 4-8: shift the inner rounded-rectangle rows up (Top only, Left/Width untouched) ---
$sh92 = Get-ShapeById $s.Shapes 92
$sh92.Top = 158.5221405029297

$sh93 = Get-ShapeById $s.Shapes 93
$sh93.Top = 187.6517333984375

$sh94 = Get-ShapeById $s.Shapes 94
$sh94.Top = 216.28843688964844

$sh95 = Get-ShapeById $s.Shapes 95
$sh95.Top = 245.24551391601562

$sh96 = Get-ShapeById $s.Shapes 96
$sh96.Top = 274.2026062011719

# --- 9: "Advanced Profiling" box grows taller + gets a new bullet ---
$sh97 = Get-ShapeById $s.Shapes 97
$sh97.TextFrame.TextRange.InsertAfter("`r• Exceptions") | Out-Null
$sh97.Top = 309.5909729003906
$sh97.Height = 58.91047286987305

# --- 10: connector between box 91 and box 2 follows the resize ---
$sh120 = Get-ShapeById $s.Shapes 120
$sh120.Top = 254.2659149169922
$sh120.Height = 0.10496063530445099
